$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")
$ws.Range("A2:L12").NumberFormat = "@"

$ws.Cells.Item(2,1).Value = "19"
$ws.Cells.Item(2,2).Value = "0"
$ws.Cells.Item(2,3).Value = "150"
$ws.Cells.Item(2,4).Value = "nan"
$ws.Cells.Item(2,5).Value = "nan"
$ws.Cells.Item(2,6).Value = "nan"
$ws.Cells.Item(2,7).Value = "nan"
$ws.Cells.Item(2,8).Value = "nan"
$ws.Cells.Item(2,9).Value = "nan"
$ws.Cells.Item(2,10).Value = "nan"
$ws.Cells.Item(2,11).Value = "nan"
$ws.Cells.Item(2,12).Value = "nan"

$ws.Cells.Item(3,1).Value = "2"
$ws.Cells.Item(3,2).Value = "151"
$ws.Cells.Item(3,3).Value = "300"
$ws.Cells.Item(3,4).Value = "nan"
$ws.Cells.Item(3,5).Value = "✔"
$ws.Cells.Item(3,6).Value = "nan"
$ws.Cells.Item(3,7).Value = "nan"
$ws.Cells.Item(3,8).Value = "nan"
$ws.Cells.Item(3,9).Value = "nan"
$ws.Cells.Item(3,10).Value = "nan"
$ws.Cells.Item(3,11).Value = "nan"
$ws.Cells.Item(3,12).Value = "5\12\2024"

$ws.Cells.Item(4,1).Value = "2"
$ws.Cells.Item(4,2).Value = "301"
$ws.Cells.Item(4,3).Value = "450"
$ws.Cells.Item(4,4).Value = "nan"
$ws.Cells.Item(4,5).Value = "nan"
$ws.Cells.Item(4,6).Value = "nan"
$ws.Cells.Item(4,7).Value = "nan"
$ws.Cells.Item(4,8).Value = "nan"
$ws.Cells.Item(4,9).Value = "nan"
$ws.Cells.Item(4,10).Value = "nan"
$ws.Cells.Item(4,11).Value = "nan"
$ws.Cells.Item(4,12).Value = "nan"

$ws.Cells.Item(5,1).Value = "2"
$ws.Cells.Item(5,2).Value = "451"
$ws.Cells.Item(5,3).Value = "550"
$ws.Cells.Item(5,4).Value = "nan"
$ws.Cells.Item(5,5).Value = "nan"
$ws.Cells.Item(5,6).Value = "nan"
$ws.Cells.Item(5,7).Value = "nan"
$ws.Cells.Item(5,8).Value = "nan"
$ws.Cells.Item(5,9).Value = "nan"
$ws.Cells.Item(5,10).Value = "nan"
$ws.Cells.Item(5,11).Value = "nan"
$ws.Cells.Item(5,12).Value = "nan"

$ws.Cells.Item(6,1).Value = "2"
$ws.Cells.Item(6,2).Value = "551"
$ws.Cells.Item(6,3).Value = "700"
$ws.Cells.Item(6,4).Value = "595"
$ws.Cells.Item(6,5).Value = "nan"
$ws.Cells.Item(6,6).Value = "nan"
$ws.Cells.Item(6,7).Value = "nan"
$ws.Cells.Item(6,8).Value = "✔"
$ws.Cells.Item(6,9).Value = "✔"
$ws.Cells.Item(6,10).Value = "nan"
$ws.Cells.Item(6,11).Value = "nan"
$ws.Cells.Item(6,12).Value = "30\4\2025"

$ws.Cells.Item(7,1).Value = "2"
$ws.Cells.Item(7,2).Value = "701"
$ws.Cells.Item(7,3).Value = "850"
$ws.Cells.Item(7,4).Value = "797"
$ws.Cells.Item(7,5).Value = "nan"
$ws.Cells.Item(7,6).Value = "✔"
$ws.Cells.Item(7,7).Value = "nan"
$ws.Cells.Item(7,8).Value = "nan"
$ws.Cells.Item(7,9).Value = "nan"
$ws.Cells.Item(7,10).Value = "nan"
$ws.Cells.Item(7,11).Value = "nan"
$ws.Cells.Item(7,12).Value = "23\8\2025"

$ws.Cells.Item(8,1).Value = "2"
$ws.Cells.Item(8,2).Value = "851"
$ws.Cells.Item(8,3).Value = "1000"
$ws.Cells.Item(8,4).Value = "nan"
$ws.Cells.Item(8,5).Value = "nan"
$ws.Cells.Item(8,6).Value = "nan"
$ws.Cells.Item(8,7).Value = "nan"
$ws.Cells.Item(8,8).Value = "nan"
$ws.Cells.Item(8,9).Value = "nan"
$ws.Cells.Item(8,10).Value = "nan"
$ws.Cells.Item(8,11).Value = "nan"
$ws.Cells.Item(8,12).Value = "nan"

$ws.Cells.Item(9,1).Value = "2"
$ws.Cells.Item(9,2).Value = "1001"
$ws.Cells.Item(9,3).Value = "1150"
$ws.Cells.Item(9,4).Value = "nan"
$ws.Cells.Item(9,5).Value = "nan"
$ws.Cells.Item(9,6).Value = "nan"
$ws.Cells.Item(9,7).Value = "nan"
$ws.Cells.Item(9,8).Value = "nan"
$ws.Cells.Item(9,9).Value = "nan"
$ws.Cells.Item(9,10).Value = "nan"
$ws.Cells.Item(9,11).Value = "nan"
$ws.Cells.Item(9,12).Value = "nan"

$ws.Cells.Item(10,1).Value = "2"
$ws.Cells.Item(10,2).Value = "1151"
$ws.Cells.Item(10,3).Value = "1300"
$ws.Cells.Item(10,4).Value = "nan"
$ws.Cells.Item(10,5).Value = "nan"
$ws.Cells.Item(10,6).Value = "nan"
$ws.Cells.Item(10,7).Value = "nan"
$ws.Cells.Item(10,8).Value = "nan"
$ws.Cells.Item(10,9).Value = "nan"
$ws.Cells.Item(10,10).Value = "nan"
$ws.Cells.Item(10,11).Value = "nan"
$ws.Cells.Item(10,12).Value = "nan"

$ws.Cells.Item(11,1).Value = "2"
$ws.Cells.Item(11,2).Value = "1301"
$ws.Cells.Item(11,3).Value = "1450"
$ws.Cells.Item(11,4).Value = "nan"
$ws.Cells.Item(11,5).Value = "nan"
$ws.Cells.Item(11,6).Value = "nan"
$ws.Cells.Item(11,7).Value = "nan"
$ws.Cells.Item(11,8).Value = "nan"
$ws.Cells.Item(11,9).Value = "nan"
$ws.Cells.Item(11,10).Value = "nan"
$ws.Cells.Item(11,11).Value = "nan"
$ws.Cells.Item(11,12).Value = "nan"

$ws.Cells.Item(12,1).Value = "2"
$ws.Cells.Item(12,2).Value = "1451"
$ws.Cells.Item(12,3).Value = "1500"
$ws.Cells.Item(12,4).Value = "nan"
$ws.Cells.Item(12,5).Value = "nan"
$ws.Cells.Item(12,6).Value = "nan"
$ws.Cells.Item(12,7).Value = "nan"
$ws.Cells.Item(12,8).Value = "nan"
$ws.Cells.Item(12,9).Value = "nan"
$ws.Cells.Item(12,10).Value = "nan"
$ws.Cells.Item(12,11).Value = "nan"
$ws.Cells.Item(12,12).Value = "nan"
